$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.3667877784032498
$ws.Cells.Item(2, 3).Value = 0.02367520258751199
$ws.Cells.Item(2, 4).Value = 0.07767967227921702
$ws.Cells.Item(2, 5).Value = 0.1459357177735754
$ws.Cells.Item(2, 7).Value = 2.132390262886673
$ws.Cells.Item(2, 8).Value = 1.679317761607024
$ws.Cells.Item(2, 11).Value = 0.31311324384248
$ws.Cells.Item(2, 13).Value = 0.240240026183244
$ws.Cells.Item(3, 2).Value = 0.3408732639094296
$ws.Cells.Item(3, 3).Value = 0.02051039519560049
$ws.Cells.Item(3, 4).Value = 0.07057444233954868
$ws.Cells.Item(3, 5).Value = 0.1341661088457755
$ws.Cells.Item(3, 7).Value = 2.048427997283454
$ws.Cells.Item(3, 8).Value = 1.642811589689131
$ws.Cells.Item(3, 11).Value = 0.2866063564770513
$ws.Cells.Item(3, 13).Value = 0.220465433772965
$ws.Cells.Item(4, 2).Value = 0.3252719707324445
$ws.Cells.Item(4, 3).Value = 0.01857294867172499
$ws.Cells.Item(4, 4).Value = 0.06625089970906117
$ws.Cells.Item(4, 5).Value = 0.1270160166652659
$ws.Cells.Item(4, 7).Value = 1.997412846125798
$ws.Cells.Item(4, 8).Value = 1.620797317247025
$ws.Cells.Item(4, 11).Value = 0.2706029756565584
$ws.Cells.Item(4, 13).Value = 0.208486487690223
$ws.Cells.Item(5, 2).Value = 0.3189920283556944
$ws.Cells.Item(5, 3).Value = 0.01778476173407739
$ws.Cells.Item(5, 4).Value = 0.06449875169529662
$ws.Cells.Item(5, 5).Value = 0.1241212093996609
$ws.Cells.Item(5, 7).Value = 1.976757533310689
$ws.Cells.Item(5, 8).Value = 1.611926454105117
$ws.Cells.Item(5, 11).Value = 0.2641494200058361
$ws.Cells.Item(5, 13).Value = 0.2036454466496522
$ws.Cells.Item(6, 2).Value = 0.3179539330385239
$ws.Cells.Item(6, 3).Value = 0.01765396153638932
$ws.Cells.Item(6, 4).Value = 0.06420839365326003
$ws.Cells.Item(6, 5).Value = 0.1236416618208409
$ws.Cells.Item(6, 7).Value = 1.973335770027063
$ws.Cells.Item(6, 8).Value = 1.610459483763606
$ws.Cells.Item(6, 11).Value = 0.2630819036488106
$ws.Cells.Item(6, 13).Value = 0.2028440273033212
$ws.Cells.Item(7, 2).Value = 0.3251869630043416
$ws.Cells.Item(7, 3).Value = 0.01856231363255034
$ws.Cells.Item(7, 4).Value = 0.06622723035981437
$ws.Cells.Item(7, 5).Value = 0.1269769001241485
$ws.Cells.Item(7, 7).Value = 1.997133741705085
$ws.Cells.Item(7, 8).Value = 1.620677277090493
$ws.Cells.Item(7, 11).Value = 0.2705156661806001
$ws.Cells.Item(7, 13).Value = 0.2084210364099803
$ws.Cells.Item(8, 2).Value = 0.3577878410586948
$ws.Cells.Item(8, 3).Value = 0.02258271477725771
$ws.Cells.Item(8, 4).Value = 0.07522161126298954
$ws.Cells.Item(8, 5).Value = 0.1418614819390811
$ws.Cells.Item(8, 7).Value = 2.103327468388699
$ws.Cells.Item(8, 8).Value = 1.66664688961157
$ws.Cells.Item(8, 11).Value = 0.3039169527487786
$ws.Cells.Item(8, 13).Value = 0.2333876964512456
$ws.Cells.Item(9, 2).Value = 0.424199050078812
$ws.Cells.Item(9, 3).Value = 0.03051773150187387
$ws.Cells.Item(9, 4).Value = 0.09317529509334577
$ws.Cells.Item(9, 5).Value = 0.1716731532543534
$ws.Cells.Item(9, 7).Value = 2.315921432497959
$ws.Cells.Item(9, 8).Value = 1.760005516646601
$ws.Cells.Item(9, 11).Value = 0.3715982677624936
$ws.Cells.Item(9, 13).Value = 0.2836613153215524
$ws.Cells.Item(10, 2).Value = 0.4745357777787262
$ws.Cells.Item(10, 3).Value = 0.03638664051641172
$ws.Cells.Item(10, 4).Value = 0.1065673761956845
$ws.Cells.Item(10, 5).Value = 0.1939801696838046
$ws.Cells.Item(10, 7).Value = 2.47489119893055
$ws.Cells.Item(10, 8).Value = 1.830608210099342
$ws.Cells.Item(10, 11).Value = 0.4226922950312542
$ws.Cells.Item(10, 13).Value = 0.3214341703158183
$ws.Cells.Item(11, 2).Value = 0.4977779857120765
$ws.Cells.Item(11, 3).Value = 0.03906680272125129
$ws.Cells.Item(11, 4).Value = 0.1127055835790571
$ws.Cells.Item(11, 5).Value = 0.2042214864370564
$ws.Cells.Item(11, 7).Value = 2.547841897648709
$ws.Cells.Item(11, 8).Value = 1.863176235010314
$ws.Cells.Item(11, 11).Value = 0.4462422845816434
$ws.Cells.Item(11, 13).Value = 0.3388077596958752
$ws.Cells.Item(12, 2).Value = 0.50662912835773
$ws.Cells.Item(12, 3).Value = 0.04008333455182367
$ws.Cells.Item(12, 4).Value = 0.1150367272933437
$ws.Cells.Item(12, 5).Value = 0.2081134874537582
$ws.Cells.Item(12, 7).Value = 2.575559626781114
$ws.Cells.Item(12, 8).Value = 1.875574505776171
$ws.Cells.Item(12, 11).Value = 0.455204816858128
$ws.Cells.Item(12, 13).Value = 0.3454146595781538
$ws.Cells.Item(13, 2).Value = 0.5047206542542142
$ws.Cells.Item(13, 3).Value = 0.03986433255639099
$ws.Cells.Item(13, 4).Value = 0.1145343723332672
$ws.Cells.Item(13, 5).Value = 0.2072746542908277
$ws.Cells.Item(13, 7).Value = 2.569585960719394
$ws.Cells.Item(13, 8).Value = 1.872901393712027
$ws.Cells.Item(13, 11).Value = 0.4532725806135147
$ws.Cells.Item(13, 13).Value = 0.3439904966232561
$ws.Cells.Item(14, 2).Value = 0.4985051737629647
$ws.Cells.Item(14, 3).Value = 0.03915040038413053
$ws.Cells.Item(14, 4).Value = 0.1128972324970476
$ws.Cells.Item(14, 5).Value = 0.2045414043395439
$ws.Cells.Item(14, 7).Value = 2.550120380986698
$ws.Cells.Item(14, 8).Value = 1.864194932217288
$ws.Cells.Item(14, 11).Value = 0.4469787394627929
$ws.Cells.Item(14, 13).Value = 0.3393507512997402
$ws.Cells.Item(15, 2).Value = 0.4947045125863099
$ws.Cells.Item(15, 3).Value = 0.03871330996524591
$ws.Cells.Item(15, 4).Value = 0.1118953178623485
$ws.Cells.Item(15, 5).Value = 0.2028690224117682
$ws.Cells.Item(15, 7).Value = 2.538209296162449
$ws.Cells.Item(15, 8).Value = 1.858870519262155
$ws.Cells.Item(15, 11).Value = 0.4431294116181164
$ws.Cells.Item(15, 13).Value = 0.3365124215373356
$ws.Cells.Item(16, 2).Value = 0.4730237954822201
$ws.Cells.Item(16, 3).Value = 0.03621170459348377
$ws.Cells.Item(16, 4).Value = 0.1061671666764141
$ws.Cells.Item(16, 5).Value = 0.1933127922801177
$ws.Cells.Item(16, 7).Value = 2.470136607492833
$ws.Cells.Item(16, 8).Value = 1.828488936747362
$ws.Cells.Item(16, 11).Value = 0.4211594698125225
$ws.Cells.Item(16, 13).Value = 0.320302641960609
$ws.Cells.Item(17, 2).Value = 0.4598117007436429
$ws.Cells.Item(17, 3).Value = 0.03467979199517401
$ws.Cells.Item(17, 4).Value = 0.1026650129315385
$ws.Cells.Item(17, 5).Value = 0.1874746168412216
$ws.Cells.Item(17, 7).Value = 2.428539901398352
$ws.Cells.Item(17, 8).Value = 1.809966689989835
$ws.Cells.Item(17, 11).Value = 0.4077606254582804
$ws.Cells.Item(17, 13).Value = 0.3104076202022128
$ws.Cells.Item(18, 2).Value = 0.4522447875399109
$ws.Cells.Item(18, 3).Value = 0.03379964470114771
$ws.Cells.Item(18, 4).Value = 0.1006549974917021
$ws.Cells.Item(18, 5).Value = 0.1841254547256526
$ws.Cells.Item(18, 7).Value = 2.404674205013293
$ws.Cells.Item(18, 8).Value = 1.79935554594033
$ws.Cells.Item(18, 11).Value = 0.4000828355744659
$ws.Cells.Item(18, 13).Value = 0.304734168425604
$ws.Cells.Item(19, 2).Value = 0.4496883025109071
$ws.Cells.Item(19, 3).Value = 0.03350180456995133
$ws.Cells.Item(19, 4).Value = 0.09997518125024385
$ws.Cells.Item(19, 5).Value = 0.1829929869344937
$ws.Cells.Item(19, 7).Value = 2.396603881993741
$ws.Cells.Item(19, 8).Value = 1.795770049920691
$ws.Cells.Item(19, 11).Value = 0.3974882111684508
$ws.Cells.Item(19, 13).Value = 0.3028162948550914
$ws.Cells.Item(20, 2).Value = 0.4612148015318382
$ws.Cells.Item(20, 3).Value = 0.03484276560925537
$ws.Cells.Item(20, 4).Value = 0.1030373741548942
$ws.Cells.Item(20, 5).Value = 0.1880951870362253
$ws.Cells.Item(20, 7).Value = 2.432961758626959
$ws.Cells.Item(20, 8).Value = 1.811934025525971
$ws.Cells.Item(20, 11).Value = 0.4091839636301415
$ws.Cells.Item(20, 13).Value = 0.3114591056542437
$ws.Cells.Item(21, 2).Value = 0.5003294557297409
$ws.Cells.Item(21, 3).Value = 0.0393600548248827
$ws.Cells.Item(21, 4).Value = 0.1133779163964164
$ws.Cells.Item(21, 5).Value = 0.2053438477218066
$ws.Cells.Item(21, 7).Value = 2.555835359701405
$ws.Cells.Item(21, 8).Value = 1.866750448772336
$ws.Cells.Item(21, 11).Value = 0.4488261769560324
$ws.Cells.Item(21, 13).Value = 0.3407127957812435
$ws.Cells.Item(22, 2).Value = 0.5261837363079565
$ws.Cells.Item(22, 3).Value = 0.04232182147421781
$ws.Cells.Item(22, 4).Value = 0.1201753842332636
$ws.Cells.Item(22, 5).Value = 0.2166976634415647
$ws.Cells.Item(22, 7).Value = 2.6366825211947
$ws.Cells.Item(22, 8).Value = 1.902958061991171
$ws.Cells.Item(22, 11).Value = 0.4749952241034521
$ws.Cells.Item(22, 13).Value = 0.3599946256903621
$ws.Cells.Item(23, 2).Value = 0.5123581002487754
$ws.Cells.Item(23, 3).Value = 0.0407401648546255
$ws.Cells.Item(23, 4).Value = 0.1165438122802556
$ws.Cells.Item(23, 5).Value = 0.2106304045214529
$ws.Cells.Item(23, 7).Value = 2.593482718267239
$ws.Cells.Item(23, 8).Value = 1.883598204465557
$ws.Cells.Item(23, 11).Value = 0.4610042987095824
$ws.Cells.Item(23, 13).Value = 0.3496884817764752
$ws.Cells.Item(24, 2).Value = 0.4605803696891542
$ws.Cells.Item(24, 3).Value = 0.03476908347411722
$ws.Cells.Item(24, 4).Value = 0.1028690190050128
$ws.Cells.Item(24, 5).Value = 0.1878146045970652
$ws.Cells.Item(24, 7).Value = 2.430962485081864
$ws.Cells.Item(24, 8).Value = 1.811044476229767
$ws.Cells.Item(24, 11).Value = 0.4085403934150804
$ws.Cells.Item(24, 13).Value = 0.310983681378751
$ws.Cells.Item(25, 2).Value = 0.4059641835170282
$ws.Cells.Item(25, 3).Value = 0.02836482077566416
$ws.Cells.Item(25, 4).Value = 0.08828361845544919
$ws.Cells.Item(25, 5).Value = 0.1635390799977046
$ws.Cells.Item(25, 7).Value = 2.257930617892811
$ws.Cells.Item(25, 8).Value = 1.734400094424387
$ws.Cells.Item(25, 11).Value = 0.3530511736886695
$ws.Cells.Item(25, 13).Value = 0.2699166198410694
